$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update L4 note text ("LED blink" -> "LED blink -this im not sure on")
$ws.Range("L4").Value = "LED blink -this im not sure on"

# C3: de-anchor the reference from $A$2 to A3 (value stays 72)
$ws.Range("C3").Formula = "=A3/B3"
$ws.Range("C3").Style = "Normal"

# Row 4: increased clock division by factor 10
$ws.Range("B4").Value = 100000
$ws.Range("C4").Formula = "=A4/B4"
$ws.Range("C4").Style = "Normal"
